# Apply "international procurement.xlsx" update:
# Replace full country names in column C with their 2-letter country codes
# (and merge "Czechia" / "Czech Republic" into a single "CZ").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("C2").Value = "CA"
    $ws.Range("C3").Value = "DK"
    $ws.Range("C4").Value = "IS"
    $ws.Range("C5").Value = "LT"
    $ws.Range("C6").Value = "NL"
    $ws.Range("C7").Value = "NO"
    $ws.Range("C8").Value = "NO"
    $ws.Range("C9").Value = "SE"
    $ws.Range("C10").Value = "UK"
    $ws.Range("C11").Value = "UK"
    $ws.Range("C12").Value = "AU"
    $ws.Range("C13").Value = "AU"
    $ws.Range("C14").Value = "BE"
    $ws.Range("C15").Value = "BE"
    $ws.Range("C16").Value = "CA"
    $ws.Range("C17").Value = "CA"
    $ws.Range("C18").Value = "HR"
    $ws.Range("C19").Value = "DK"
    $ws.Range("C20").Value = "DE"
    $ws.Range("C21").Value = "IS"
    $ws.Range("C22").Value = "IS"
    $ws.Range("C23").Value = "JP"
    $ws.Range("C24").Value = "LU"
    $ws.Range("C25").Value = "LU"
    $ws.Range("C26").Value = "NL"
    $ws.Range("C27").Value = "NZ"
    $ws.Range("C28").Value = "NZ"
    $ws.Range("C29").Value = "NZ"
    $ws.Range("C30").Value = "NZ"
    $ws.Range("C31").Value = "NO"
    $ws.Range("C32").Value = "RO"
    $ws.Range("C33").Value = "RO"
    $ws.Range("C34").Value = "SE"
    $ws.Range("C35").Value = "SE"
    $ws.Range("C36").Value = "UK"
    $ws.Range("C37").Value = "UK"
    $ws.Range("C38").Value = "DK"
    $ws.Range("C39").Value = "SE"
    $ws.Range("C40").Value = "FI"
    $ws.Range("C41").Value = "NO"
    $ws.Range("C42").Value = "AT"
    $ws.Range("C43").Value = "BE"
    $ws.Range("C44").Value = "BG"
    $ws.Range("C45").Value = "HR"
    $ws.Range("C46").Value = "CY"
    $ws.Range("C47").Value = "CZ"
    $ws.Range("C48").Value = "DK"
    $ws.Range("C49").Value = "EE"
    $ws.Range("C50").Value = "FI"
    $ws.Range("C51").Value = "FR"
    $ws.Range("C52").Value = "DE"
    $ws.Range("C53").Value = "GR"
    $ws.Range("C54").Value = "HU"
    $ws.Range("C55").Value = "IE"
    $ws.Range("C56").Value = "IT"
    $ws.Range("C57").Value = "LV"
    $ws.Range("C58").Value = "LT"
    $ws.Range("C59").Value = "LU"
    $ws.Range("C60").Value = "MT"
    $ws.Range("C61").Value = "NL"
    $ws.Range("C62").Value = "PL"
    $ws.Range("C63").Value = "PT"
    $ws.Range("C64").Value = "RO"
    $ws.Range("C65").Value = "SK"
    $ws.Range("C66").Value = "SL"
    $ws.Range("C67").Value = "ES"
    $ws.Range("C68").Value = "SE"
    $ws.Range("C69").Value = "NL"
    $ws.Range("C70").Value = "DK"
    $ws.Range("C71").Value = "DK"
    $ws.Range("C72").Value = "NL"
    $ws.Range("C73").Value = "DK"
    $ws.Range("C74").Value = "NL"
    $ws.Range("C75").Value = "CZ"

# Reposition the view/selection to match the saved state of the file
# (top-left cell scrolled down, active cell on the last data row).
$ws.Application.ActiveWindow.ScrollRow = 58
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C75").Select()
